# Final set of survival and size data
# Fill in day19 (W) / day20 (X) observations for rows 29-36 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row -> (W value, X value)
$values = @{
    29 = @(1, 1)
    30 = @(1, 1)
    31 = @(1, 1)
    32 = @(0, 0)
    33 = @(0, 0)
    34 = @(1, 1)
    35 = @(0, 0)
    36 = @(0, 0)
}

foreach ($row in $values.Keys | Sort-Object) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 23).Value = $pair[0]   # column W = day19
    $ws.Cells.Item($row, 24).Value = $pair[1]   # column X = day20
}

# Reflect the updated view state (scroll position + active selection) recorded in the commit.
$win = $excel.ActiveWindow
$win.ScrollColumn = 15   # topLeftCell -> O1
$win.ScrollRow = 1
$ws.Range("Y8").Select()
